# Scheduled market-data refresh: update the H-N (price/profit) columns for a
# handful of Leve rows across several sheets with freshly scraped values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 399.70587
$ws.Range("I12").Value = 342.16666
$ws.Range("K12").Value = 342.16666
$ws.Range("M12").Value = -172.16666
$ws.Range("H62").Value = 2249.7
$ws.Range("I62").Value = 1960.4
$ws.Range("K62").Value = 1960.4
$ws.Range("M62").Value = -1336.4
$ws.Range("H65").Value = 2249.7
$ws.Range("I65").Value = 1960.4
$ws.Range("K65").Value = 9802
$ws.Range("M65").Value = -6682
$ws.Range("H68").Value = 169999
$ws.Range("J68").Value = 169999
$ws.Range("L68").Value = 169999
$ws.Range("N68").Value = -171497
$ws.Range("H70").Value = 73064.21000000001
$ws.Range("I70").Value = 1777.7778
$ws.Range("J70").Value = 201379.8
$ws.Range("K70").Value = 5333.3334
$ws.Range("L70").Value = 604139.3999999999
$ws.Range("M70").Value = -5063.3334
$ws.Range("N70").Value = -604679.3999999999
$ws.Range("H71").Value = 169999
$ws.Range("J71").Value = 169999
$ws.Range("L71").Value = 509997
$ws.Range("N71").Value = -517485
$ws.Range("H73").Value = 73064.21000000001
$ws.Range("I73").Value = 1777.7778
$ws.Range("J73").Value = 201379.8
$ws.Range("K73").Value = 5333.3334
$ws.Range("L73").Value = 604139.3999999999
$ws.Range("M73").Value = -4397.3334
$ws.Range("N73").Value = -606011.3999999999
$ws.Range("H138").Value = 6695.5366
$ws.Range("J138").Value = 6791.5
$ws.Range("L138").Value = 20374.5
$ws.Range("N138").Value = -30654.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2392.375
$ws.Range("I2").Value = 2023.9062
$ws.Range("J2").Value = 3866.25
$ws.Range("K2").Value = 2023.9062
$ws.Range("L2").Value = 3866.25
$ws.Range("M2").Value = -1910.9062
$ws.Range("N2").Value = -4092.25
$ws.Range("H45").Value = 2757.9
$ws.Range("I45").Value = 2026.7273
$ws.Range("K45").Value = 2026.7273
$ws.Range("M45").Value = -1649.7273
$ws.Range("H61").Value = 4730.609
$ws.Range("I61").Value = 4658.2104
$ws.Range("K61").Value = 4658.2104
$ws.Range("M61").Value = -4446.2104
$ws.Range("H116").Value = 2392.375
$ws.Range("I116").Value = 2023.9062
$ws.Range("J116").Value = 3866.25
$ws.Range("K116").Value = 2023.9062
$ws.Range("L116").Value = 3866.25
$ws.Range("M116").Value = 270.0938000000001
$ws.Range("N116").Value = -8454.25
$ws.Range("H136").Value = 4730.609
$ws.Range("I136").Value = 4658.2104
$ws.Range("K136").Value = 13974.6312
$ws.Range("M136").Value = -11424.6312

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2392.375
$ws.Range("I3").Value = 2023.9062
$ws.Range("J3").Value = 3866.25
$ws.Range("K3").Value = 2023.9062
$ws.Range("L3").Value = 3866.25
$ws.Range("M3").Value = -1909.9062
$ws.Range("N3").Value = -4094.25
$ws.Range("H105").Value = 2766.4
$ws.Range("I105").Value = 2639.5557
$ws.Range("J105").Value = 3908
$ws.Range("K105").Value = 2639.5557
$ws.Range("L105").Value = 3908
$ws.Range("M105").Value = -892.5556999999999
$ws.Range("N105").Value = -7402
$ws.Range("H134").Value = 1984.1428
$ws.Range("I134").Value = 1829.1666
$ws.Range("J134").Value = 2914
$ws.Range("K134").Value = 5487.4998
$ws.Range("L134").Value = 8742
$ws.Range("M134").Value = -2952.4998
$ws.Range("N134").Value = -13812

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5443.12
$ws.Range("I31").Value = 5422.276
$ws.Range("J31").Value = 5471.905
$ws.Range("K31").Value = 5422.276
$ws.Range("L31").Value = 5471.905
$ws.Range("M31").Value = -5127.276
$ws.Range("N31").Value = -6061.905
$ws.Range("H34").Value = 5443.12
$ws.Range("I34").Value = 5422.276
$ws.Range("J34").Value = 5471.905
$ws.Range("K34").Value = 5422.276
$ws.Range("L34").Value = 5471.905
$ws.Range("M34").Value = -5220.276
$ws.Range("N34").Value = -5875.905
$ws.Range("H86").Value = 7667.7144
$ws.Range("I86").Value = 5292.3335
$ws.Range("J86").Value = 9449.25
$ws.Range("K86").Value = 5292.3335
$ws.Range("L86").Value = 9449.25
$ws.Range("M86").Value = -4169.3335
$ws.Range("N86").Value = -11695.25
$ws.Range("H89").Value = 7667.7144
$ws.Range("I89").Value = 5292.3335
$ws.Range("J89").Value = 9449.25
$ws.Range("K89").Value = 26461.6675
$ws.Range("L89").Value = 47246.25
$ws.Range("M89").Value = -20845.6675
$ws.Range("N89").Value = -58478.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 555.5
$ws.Range("J9").Value = 650
$ws.Range("L9").Value = 650
$ws.Range("N9").Value = -990
$ws.Range("H80").Value = 5360
$ws.Range("I80").Value = 3234.6667
$ws.Range("J80").Value = 6776.8887
$ws.Range("K80").Value = 3234.6667
$ws.Range("L80").Value = 6776.8887
$ws.Range("M80").Value = -2236.6667
$ws.Range("N80").Value = -8772.8887
$ws.Range("H83").Value = 5360
$ws.Range("I83").Value = 3234.6667
$ws.Range("J83").Value = 6776.8887
$ws.Range("K83").Value = 16173.3335
$ws.Range("L83").Value = 33884.4435
$ws.Range("M83").Value = -11181.3335
$ws.Range("N83").Value = -43868.4435
$ws.Range("H113").Value = 12026.786
$ws.Range("I113").Value = 6373.5
$ws.Range("K113").Value = 6373.5
$ws.Range("M113").Value = -4203.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 354.61905
$ws.Range("I55").Value = 242.75
$ws.Range("J55").Value = 503.77777
$ws.Range("K55").Value = 242.75
$ws.Range("L55").Value = 503.77777
$ws.Range("M55").Value = -69.75
$ws.Range("N55").Value = -849.7777699999999
$ws.Range("H100").Value = 3496.4167
$ws.Range("I100").Value = 1658.8334
$ws.Range("K100").Value = 1658.8334
$ws.Range("M100").Value = -1117.8334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 5668
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 5668
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 5668
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -5898
$ws.Range("H63").Value = 10999
$ws.Range("J63").Value = 10999
$ws.Range("L63").Value = 10999
$ws.Range("N63").Value = -12247
$ws.Range("H64").Value = 49998.5
$ws.Range("J64").Value = 49998.5
$ws.Range("L64").Value = 49998.5
$ws.Range("N64").Value = -50494.5
$ws.Range("H66").Value = 10999
$ws.Range("J66").Value = 10999
$ws.Range("L66").Value = 32997
$ws.Range("N66").Value = -39237
$ws.Range("H67").Value = 49998.5
$ws.Range("J67").Value = 49998.5
$ws.Range("L67").Value = 49998.5
$ws.Range("N67").Value = -51714.5
$ws.Range("H68").Value = 37635
$ws.Range("J68").Value = 37635
$ws.Range("L68").Value = 37635
$ws.Range("N68").Value = -39257
$ws.Range("H71").Value = 37635
$ws.Range("J71").Value = 37635
$ws.Range("L71").Value = 112905
$ws.Range("N71").Value = -121017
$ws.Range("H126").Value = 2194.0588
$ws.Range("I126").Value = 1906.6
$ws.Range("K126").Value = 5719.799999999999
$ws.Range("M126").Value = -3249.799999999999
$ws.Range("H132").Value = 2299.6086
$ws.Range("I132").Value = 2107.4055
$ws.Range("J132").Value = 3089.7778
$ws.Range("K132").Value = 6322.2165
$ws.Range("L132").Value = 9269.3334
$ws.Range("M132").Value = -3792.2165
$ws.Range("N132").Value = -14329.3334
